$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Remove the columns that are no longer used:
# Project (B), ID (C), PCT (E), P_STATUS (F), RO Number (N)
# Deleting from right to left keeps earlier column letters valid.
$ws.Columns("N").Delete()
$ws.Columns("F").Delete()
$ws.Columns("E").Delete()
$ws.Columns("C").Delete()
$ws.Columns("B").Delete()

# After the deletions the remaining columns are:
# A Proposal, B Faculty, C Sponsor, D Allocated Amt, E Total Cost, F Funded?,
# G Long Descr, H Begin Date, I End Date, J Submit Date, K Principal Investigators

# Rename the header fields to their new labels.
$ws.Range("A1").Value = "Proposal_ID"
$ws.Range("G1").Value = "Title"

# Make the new Proposal_ID header bold.
$ws.Range("A1").Font.Bold = $true

# Adjust column widths for the new layout.
$ws.Columns("A").ColumnWidth = 9.833333333333333
$ws.Columns("C").ColumnWidth = 13.5
$ws.Columns("D").ColumnWidth = 13.5
$ws.Columns("E").ColumnWidth = 9.5
$ws.Columns("J").ColumnWidth = 11.5

# Make the Data sheet the active sheet/tab, with C1 selected.
$ws.Activate()
$ws.Range("C1").Select()
